$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting cell with the new commit text
$ws.Range("E8").Value = "GIT UPDATE"

# Select E8 as the active cell, matching saved view state
$ws.Range("E8").Select()
